$d = $word.ActiveDocument

# Locate the run that reads "Website processing error" (red, bold) so we can
# append a new run right after it in the same paragraph.
$r = $d.Content
$found = $r.Find.Execute("Website processing error", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Build a minimal WordprocessingML package fragment describing the new
    # run: bold, bold-complex-script and red-colored text "-> resolved".
    # Calling InsertXML on the (non-collapsed) Find result range appends the
    # new run immediately after the matched text, inside the same paragraph,
    # leaving the matched run itself untouched.
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' +
           '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/></w:rPr><w:t>-&gt; resolved</w:t></w:r>' +
           '</w:p></w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $r.InsertXML($xml)
}
